# edit.ps1
# Applies updated "想去人数" (interest count, column F) values to the
# "展览" (sheet1), "演出" (sheet2) and "全部类型" (sheet4) worksheets,
# matching the regenerated data snapshot (commit: "Update gh-pages to output generated at 456a3b4").
# "本地生活" (sheet3) has no data rows, so no changes are needed there.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Cells.Item(2, 6).Value = 4739  # row 2: 4727 -> 4739
$ws.Cells.Item(3, 6).Value = 2739  # row 3: 2733 -> 2739
$ws.Cells.Item(5, 6).Value = 2777  # row 5: 2766 -> 2777
$ws.Cells.Item(9, 6).Value = 1725  # row 9: 1724 -> 1725
$ws.Cells.Item(10, 6).Value = 748  # row 10: 747 -> 748
$ws.Cells.Item(12, 6).Value = 231  # row 12: 227 -> 231
$ws.Cells.Item(13, 6).Value = 404  # row 13: 402 -> 404
$ws.Cells.Item(14, 6).Value = 1060  # row 14: 1059 -> 1060
$ws.Cells.Item(15, 6).Value = 297  # row 15: 296 -> 297
$ws.Cells.Item(17, 6).Value = 70  # row 17: 68 -> 70
$ws.Cells.Item(18, 6).Value = 529  # row 18: 528 -> 529
$ws.Cells.Item(19, 6).Value = 529  # row 19: 528 -> 529
$ws.Cells.Item(22, 6).Value = 653  # row 22: 649 -> 653
$ws.Cells.Item(23, 6).Value = 737  # row 23: 736 -> 737
$ws.Cells.Item(24, 6).Value = 147  # row 24: 145 -> 147
$ws.Cells.Item(26, 6).Value = 508  # row 26: 506 -> 508
$ws.Cells.Item(27, 6).Value = 13  # row 27: 10 -> 13
$ws.Cells.Item(29, 6).Value = 1513  # row 29: 1493 -> 1513
$ws.Cells.Item(30, 6).Value = 329  # row 30: 326 -> 329
$ws.Cells.Item(32, 6).Value = 1461  # row 32: 1450 -> 1461
$ws.Cells.Item(33, 6).Value = 139  # row 33: 111 -> 139
$ws.Cells.Item(34, 6).Value = 2323  # row 34: 2311 -> 2323
$ws.Cells.Item(35, 6).Value = 385  # row 35: 383 -> 385
$ws.Cells.Item(37, 6).Value = 605  # row 37: 604 -> 605
$ws.Cells.Item(38, 6).Value = 114  # row 38: 113 -> 114
$ws.Cells.Item(41, 6).Value = 778  # row 41: 776 -> 778
$ws.Cells.Item(42, 6).Value = 1477  # row 42: 1474 -> 1477
$ws.Cells.Item(43, 6).Value = 205  # row 43: 203 -> 205
$ws.Cells.Item(45, 6).Value = 491  # row 45: 488 -> 491
$ws.Cells.Item(46, 6).Value = 34  # row 46: 30 -> 34
$ws.Cells.Item(47, 6).Value = 76  # row 47: 77 -> 76
$ws.Cells.Item(48, 6).Value = 107  # row 48: 106 -> 107

$ws = $wb.Worksheets.Item("演出")
$ws.Cells.Item(12, 6).Value = 35  # row 12: 34 -> 35

$ws = $wb.Worksheets.Item("全部类型")
$ws.Cells.Item(2, 6).Value = 4739  # row 2: 4727 -> 4739
$ws.Cells.Item(3, 6).Value = 2739  # row 3: 2733 -> 2739
$ws.Cells.Item(4, 6).Value = 2777  # row 4: 2766 -> 2777
$ws.Cells.Item(5, 6).Value = 1725  # row 5: 1724 -> 1725
$ws.Cells.Item(8, 6).Value = 748  # row 8: 747 -> 748
$ws.Cells.Item(10, 6).Value = 231  # row 10: 227 -> 231
$ws.Cells.Item(11, 6).Value = 404  # row 11: 402 -> 404
$ws.Cells.Item(12, 6).Value = 1060  # row 12: 1059 -> 1060
$ws.Cells.Item(13, 6).Value = 297  # row 13: 296 -> 297
$ws.Cells.Item(15, 6).Value = 70  # row 15: 68 -> 70
$ws.Cells.Item(16, 6).Value = 529  # row 16: 528 -> 529
$ws.Cells.Item(17, 6).Value = 529  # row 17: 528 -> 529
$ws.Cells.Item(19, 6).Value = 653  # row 19: 649 -> 653
$ws.Cells.Item(20, 6).Value = 737  # row 20: 736 -> 737
$ws.Cells.Item(21, 6).Value = 147  # row 21: 145 -> 147
$ws.Cells.Item(26, 6).Value = 508  # row 26: 506 -> 508
$ws.Cells.Item(28, 6).Value = 1513  # row 28: 1494 -> 1513
$ws.Cells.Item(29, 6).Value = 329  # row 29: 326 -> 329
$ws.Cells.Item(33, 6).Value = 2323  # row 33: 2311 -> 2323
$ws.Cells.Item(34, 6).Value = 385  # row 34: 383 -> 385
$ws.Cells.Item(38, 6).Value = 35  # row 38: 34 -> 35
$ws.Cells.Item(39, 6).Value = 605  # row 39: 604 -> 605
$ws.Cells.Item(40, 6).Value = 114  # row 40: 113 -> 114
$ws.Cells.Item(43, 6).Value = 778  # row 43: 776 -> 778
$ws.Cells.Item(44, 6).Value = 1477  # row 44: 1474 -> 1477
$ws.Cells.Item(46, 6).Value = 205  # row 46: 203 -> 205
$ws.Cells.Item(47, 6).Value = 491  # row 47: 488 -> 491
$ws.Cells.Item(48, 6).Value = 76  # row 48: 77 -> 76
$ws.Cells.Item(49, 6).Value = 107  # row 49: 106 -> 107

$wb.Save()
